# Auto-generated: apply scheduled market-data refresh to Sheets (Siren_Profits workbook)
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H-N) for specific leve rows
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 370.9091
$ws.Range("I41").Value = 207.25
$ws.Range("K41").Value = 207.25
$ws.Range("M41").Value = 232.75

$ws.Range("H62").Value = 19768.428
$ws.Range("I62").Value = 3683.6667
$ws.Range("J62").Value = 31832
$ws.Range("K62").Value = 3683.6667
$ws.Range("L62").Value = 31832
$ws.Range("M62").Value = -3059.6667
$ws.Range("N62").Value = -33080

$ws.Range("H65").Value = 19768.428
$ws.Range("I65").Value = 3683.6667
$ws.Range("J65").Value = 31832
$ws.Range("K65").Value = 18418.3335
$ws.Range("L65").Value = 159160
$ws.Range("M65").Value = -15298.3335
$ws.Range("N65").Value = -165400

$ws.Range("H129").Value = 1748.9231
$ws.Range("J129").Value = 4550
$ws.Range("L129").Value = 13650
$ws.Range("N129").Value = -23650

$ws.Range("H132").Value = 3353.1365
$ws.Range("I132").Value = 3174.9736
$ws.Range("J132").Value = 4481.5
$ws.Range("K132").Value = 9524.9208
$ws.Range("L132").Value = 13444.5
$ws.Range("M132").Value = -6994.9208
$ws.Range("N132").Value = -18504.5

$ws.Range("H135").Value = 2575.12
$ws.Range("I135").Value = 2549.6191
$ws.Range("J135").Value = 2709
$ws.Range("K135").Value = 22946.5719
$ws.Range("L135").Value = 24381
$ws.Range("M135").Value = -20411.5719
$ws.Range("N135").Value = -29451

$ws.Range("H138").Value = 4268.793
$ws.Range("I138").Value = 1756.4615
$ws.Range("J138").Value = 4994.5776
$ws.Range("K138").Value = 5269.3845
$ws.Range("L138").Value = 14983.7328
$ws.Range("M138").Value = -129.3845000000001
$ws.Range("N138").Value = -25263.7328

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3106.75
$ws.Range("I2").Value = 2517.4
$ws.Range("J2").Value = 3527.7144
$ws.Range("K2").Value = 2517.4
$ws.Range("L2").Value = 3527.7144
$ws.Range("M2").Value = -2404.4
$ws.Range("N2").Value = -3753.7144

$ws.Range("H32").Value = 3987.4922
$ws.Range("I32").Value = 3551.4355
$ws.Range("K32").Value = 3551.4355
$ws.Range("M32").Value = -3264.4355

$ws.Range("H102").Value = 13060.579
$ws.Range("I102").Value = 15543.8
$ws.Range("K102").Value = 15543.8
$ws.Range("M102").Value = -13921.8

$ws.Range("H116").Value = 3106.75
$ws.Range("I116").Value = 2517.4
$ws.Range("J116").Value = 3527.7144
$ws.Range("K116").Value = 2517.4
$ws.Range("L116").Value = 3527.7144
$ws.Range("M116").Value = -223.4000000000001
$ws.Range("N116").Value = -8115.7144

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3106.75
$ws.Range("I3").Value = 2517.4
$ws.Range("J3").Value = 3527.7144
$ws.Range("K3").Value = 2517.4
$ws.Range("L3").Value = 3527.7144
$ws.Range("M3").Value = -2403.4
$ws.Range("N3").Value = -3755.7144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4393.7856
$ws.Range("J31").Value = 4750
$ws.Range("L31").Value = 4750
$ws.Range("N31").Value = -5340

$ws.Range("H34").Value = 4393.7856
$ws.Range("J34").Value = 4750
$ws.Range("L34").Value = 4750
$ws.Range("N34").Value = -5154

$ws.Range("H58").Value = 2034.3103
$ws.Range("I58").Value = 1796.3478
$ws.Range("J58").Value = 2946.5
$ws.Range("K58").Value = 1796.3478
$ws.Range("L58").Value = 2946.5
$ws.Range("M58").Value = -1593.3478
$ws.Range("N58").Value = -3352.5

$ws.Range("H86").Value = 10564.8
$ws.Range("J86").Value = 17966.334
$ws.Range("L86").Value = 17966.334
$ws.Range("N86").Value = -20212.334

$ws.Range("H89").Value = 10564.8
$ws.Range("J89").Value = 17966.334
$ws.Range("L89").Value = 89831.67
$ws.Range("N89").Value = -101063.67

$ws.Range("H99").Value = 14518482
$ws.Range("J99").Value = 5999.3335
$ws.Range("L99").Value = 5999.3335
$ws.Range("N99").Value = -8995.333500000001

$ws.Range("H126").Value = 14518482
$ws.Range("J126").Value = 5999.3335
$ws.Range("L126").Value = 17998.0005
$ws.Range("N126").Value = -22938.0005

$ws.Range("H132").Value = 10511.021
$ws.Range("I132").Value = 1274.775
$ws.Range("K132").Value = 3824.325
$ws.Range("M132").Value = -1294.325

$ws.Range("H136").Value = 2034.3103
$ws.Range("I136").Value = 1796.3478
$ws.Range("J136").Value = 2946.5
$ws.Range("K136").Value = 5389.0434
$ws.Range("L136").Value = 8839.5
$ws.Range("M136").Value = -2839.0434
$ws.Range("N136").Value = -13939.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1041.6666
$ws.Range("I11").Value = 1150
$ws.Range("J11").Value = 987.5
$ws.Range("K11").Value = 3450
$ws.Range("L11").Value = 2962.5
$ws.Range("M11").Value = -3310
$ws.Range("N11").Value = -3242.5

$ws.Range("H38").Value = 1522.44
$ws.Range("I38").Value = 334
$ws.Range("K38").Value = 1002
$ws.Range("M38").Value = -655

$ws.Range("H59").Value = 2891
$ws.Range("J59").Value = 5052.5
$ws.Range("L59").Value = 15157.5
$ws.Range("N59").Value = -16237.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 30000
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H59").Value = 7775
$ws.Range("J59").Value = 6550
$ws.Range("L59").Value = 6550
$ws.Range("N59").Value = -7716

$ws.Range("H107").Value = 708.5
$ws.Range("I107").Value = 696.7143
$ws.Range("J107").Value = 725
$ws.Range("K107").Value = 696.7143
$ws.Range("L107").Value = 725
$ws.Range("M107").Value = 1223.2857
$ws.Range("N107").Value = -4565

$ws.Range("H113").Value = 13368.1
$ws.Range("J113").Value = 2343
$ws.Range("L113").Value = 2343
$ws.Range("N113").Value = -6683

$ws.Range("H122").Value = 16077.167
$ws.Range("I122").Value = 10642.6
$ws.Range("J122").Value = 43250
$ws.Range("K122").Value = 31927.8
$ws.Range("L122").Value = 129750
$ws.Range("M122").Value = -29477.8
$ws.Range("N122").Value = -134650

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 7026.2607
$ws.Range("I61").Value = 6022.3335
$ws.Range("J61").Value = 8121.4546
$ws.Range("K61").Value = 6022.3335
$ws.Range("L61").Value = 8121.4546
$ws.Range("M61").Value = -5820.3335
$ws.Range("N61").Value = -8525.454600000001

$ws.Range("H113").Value = 7026.2607
$ws.Range("I113").Value = 6022.3335
$ws.Range("J113").Value = 8121.4546
$ws.Range("K113").Value = 6022.3335
$ws.Range("L113").Value = 8121.4546
$ws.Range("M113").Value = -3852.3335
$ws.Range("N113").Value = -12461.4546

$ws.Range("H132").Value = 467559.12
$ws.Range("I132").Value = 498130.1
$ws.Range("J132").Value = 8994.5
$ws.Range("K132").Value = 1494390.3
$ws.Range("L132").Value = 26983.5
$ws.Range("M132").Value = -1491860.3
$ws.Range("N132").Value = -32043.5

$ws.Range("H136").Value = 6212.3184
$ws.Range("J136").Value = 15827
$ws.Range("L136").Value = 47481
$ws.Range("N136").Value = -52581

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 6379.75
$ws.Range("I33").Value = 6379.75
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 6379.75
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -6129.75
$ws.Range("N33").ClearContents()

$ws.Range("H36").Value = 6379.75
$ws.Range("I36").Value = 6379.75
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 6379.75
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -6129.75
$ws.Range("N36").ClearContents()

$ws.Range("H100").Value = 42139.867
$ws.Range("I100").Value = 31319.8
$ws.Range("J100").Value = 63780
$ws.Range("K100").Value = 62639.6
$ws.Range("L100").Value = 127560
$ws.Range("M100").Value = -62098.6
$ws.Range("N100").Value = -128642

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws.Range("H135").Value = 137081.86
$ws.Range("J135").Value = 137081.86
$ws.Range("L135").Value = 137081.86
$ws.Range("N135").Value = -147221.86

$ws.Range("H136").Value = 2597.9546
$ws.Range("I136").Value = 1786.7222
$ws.Range("K136").Value = 5360.1666
$ws.Range("M136").Value = -2810.1666
